$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "ParentsEmail"
$ws.Range("M2").Select()
